$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text (as in source)
$textCells = @("D5", "D6", "D7", "D12", "D14", "D20", "D23", "D24", "D25", "D29", "D31", "D32", "D35", "D37", "D39", "D43", "D47", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.868.04"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "3.322.57"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "578.88"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "174.85"
$ws.Range("E6").Value = "  -4.70%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "3.318.39"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "45.39"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "659.83"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "3.860.87"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "67.777.62"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "3.327.34"
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "5.34"
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").Value = "16.86"
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("D25").Value = "98.10"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  -4.16%  "
$ws.Range("E27").Value = "  -4.56%  "
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("D29").Value = "33.28"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").Value = "7.22"
$ws.Range("E31").Value = "  +7.50%  "
$ws.Range("D32").Value = "566.49"
$ws.Range("E32").Value = "  -6.09%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").Value = "3.680.92"
$ws.Range("E36").Value = "  -7.18%  "
$ws.Range("D37").Value = "56.30"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  -8.62%  "
$ws.Range("D39").Value = "34.34"
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").Value = "3.30"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "0.0₃0660"
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.68"
$ws.Range("E51").Value = "  +9.13%  "

# Reset style back to default (no explicit style) for cells where we forced text format
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
